$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: "This is a Microsoft word document." gains two
#    trailing spaces, then three new runs colored C00000 are appended:
#       "(This is a change – Ve"
#       "rsion for branch alternate"
#       ")"
# ---------------------------------------------------------------------------

$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End
# Range covering just the paragraph's text, excluding the trailing pilcrow.
$textRng = $d.Range($p1Start, $p1End - 1)
[void]$textRng.InsertAfter("  ")

# Position right after the two trailing spaces we just appended (still
# before the paragraph mark) — this is where the colored runs get inserted.
$insertPos = $p1Start + 34 + 2

$run1Text = "(This is a change " + [string][char]0x2013 + " Ve"
$r1 = $d.Range($insertPos, $insertPos)
[void]$r1.InsertAfter($run1Text)
$r1.Font.Color = 192

$pos2 = $insertPos + $run1Text.Length
$run2Text = "rsion for branch alternate"
$r2 = $d.Range($pos2, $pos2)
[void]$r2.InsertAfter($run2Text)
$r2.Font.Color = 192

$pos3 = $pos2 + $run2Text.Length
$run3Text = ")"
$r3 = $d.Range($pos3, $pos3)
[void]$r3.InsertAfter($run3Text)
$r3.Font.Color = 192

# ---------------------------------------------------------------------------
# 2. Append a new, bare paragraph (just light-grey shading, no text/style
#    inheritance) right after the last paragraph in the document.
# ---------------------------------------------------------------------------

$endRng = $d.Content
$endRng.Collapse(0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
[void]$endRng.InsertXML($newParaXml)
